$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
$ws.Range("A8").Value = "Volume 31   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/10/2024  Through  6/16/2024"

# --- Special cells that switch between numeric and text (shared-string) representation ---
# Copy format+content from a donor cell that already has the desired style/type, then set numeric value where needed
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))
$ws.Range("C18").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 2
$ws.Range("H18").Copy($ws.Range("E18"))
$ws.Range("E18").Value = -50

# --- Bulk numeric updates ---
$ws.Range("N14").Value = -95.238095238095
$ws.Range("G15").Value = 5
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = -18.181818181818
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = -84.482758620689
$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 700
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 183.333333333333
$ws.Range("I16").Value = 86
$ws.Range("J16").Value = 69
$ws.Range("K16").Value = 24.637681159420
$ws.Range("L16").Value = -14.851485148514
$ws.Range("M16").Value = -50.574712643678
$ws.Range("N16").Value = -91.041666666666
$ws.Range("C17").Value = 21
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 162.5
$ws.Range("F17").Value = 46
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 91.666666666666
$ws.Range("I17").Value = 147
$ws.Range("J17").Value = 148
$ws.Range("K17").Value = -0.675675675675
$ws.Range("L17").Value = -15.028901734104
$ws.Range("M17").Value = -10.909090909090
$ws.Range("N17").Value = -62.404092071611
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 125
$ws.Range("I18").Value = 59
$ws.Range("J18").Value = 68
$ws.Range("K18").Value = -13.235294117647
$ws.Range("L18").Value = -21.333333333333
$ws.Range("M18").Value = -52.8
$ws.Range("N18").Value = -95.485845447589
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -61.538461538461
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = -4.545454545454
$ws.Range("I19").Value = 239
$ws.Range("J19").Value = 263
$ws.Range("K19").Value = -9.125475285171
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -14.336917562724
$ws.Range("N19").Value = -52.766798418972
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -8.333333333333
$ws.Range("I20").Value = 65
$ws.Range("J20").Value = 47
$ws.Range("K20").Value = 38.297872340425
$ws.Range("L20").Value = 12.068965517241
$ws.Range("M20").Value = -30.851063829787
$ws.Range("N20").Value = -94.519392917369
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 40.740740740740
$ws.Range("F21").Value = 126
$ws.Range("G21").Value = 95
$ws.Range("H21").Value = 32.631578947368
$ws.Range("I21").Value = 606
$ws.Range("J21").Value = 609
$ws.Range("K21").Value = -0.492610837438
$ws.Range("L21").Value = -8.872180451127
$ws.Range("M21").Value = -28.705882352941
$ws.Range("N21").Value = -86.317453149695
$ws.Range("G22").Value = 2
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = -10
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 53
$ws.Range("E24").Value = -47.169811320754
$ws.Range("F24").Value = 126
$ws.Range("G24").Value = 177
$ws.Range("H24").Value = -28.813559322033
$ws.Range("I24").Value = 800
$ws.Range("J24").Value = 838
$ws.Range("K24").Value = -4.534606205250
$ws.Range("L24").Value = 37.457044673539
$ws.Range("M24").Value = 42.857142857142
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 34
$ws.Range("E25").Value = -52.941176470588
$ws.Range("F25").Value = 82
$ws.Range("G25").Value = 93
$ws.Range("H25").Value = -11.827956989247
$ws.Range("I25").Value = 498
$ws.Range("J25").Value = 482
$ws.Range("K25").Value = 3.319502074688
$ws.Range("L25").Value = 102.439024390244
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 60
$ws.Range("G26").Value = 48
$ws.Range("H26").Value = 25
$ws.Range("I26").Value = 304
$ws.Range("J26").Value = 276
$ws.Range("K26").Value = 10.144927536231
$ws.Range("L26").Value = 17.829457364341
$ws.Range("M26").Value = -11.370262390670
$ws.Range("G27").Value = 5
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = -25
$ws.Range("L27").Value = -44.444444444444
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -20
$ws.Range("I28").Value = 34
$ws.Range("J28").Value = 27
$ws.Range("K28").Value = 25.925925925925
$ws.Range("L28").Value = 9.677419354838
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = -57.142857142857
$ws.Range("M29").Value = -85.714285714285
$ws.Range("N29").Value = -94.915254237288
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -33.333333333333
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -50
$ws.Range("M30").Value = -82.352941176470
$ws.Range("N30").Value = -93.617021276595
$ws.Range("C31").Value = 1
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 300
$ws.Range("I31").Value = 9
$ws.Range("K31").Value = 12.5
$ws.Range("L31").Value = 28.571428571428
